# Apply capital structure database update for Jamaica Insurance (Prop/Cas.) sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row identities rotate: row3->Key, row4->General Accident, row5->Ironrock
$ws.Range("B3").Value = "Key Insurance Company Limited (JMSE:KEY)"
$ws.Range("B4").Value = "General Accident Insurance Company Jamaica Limited (JMSE:GENAC)"
$ws.Range("B5").Value = "Ironrock Insurance Company Limited (JMSE:ROC)"

# Cells with no data in the refreshed row (previously populated)
$ws.Range("D3").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("T3").ClearContents()

# Refreshed metric values
$ws.Range("D2").Value = 0.189
$ws.Range("E2").Value = 0.214
$ws.Range("G2").Value = -0.0248
$ws.Range("H2").Value = -0.0248
$ws.Range("I2").Value = -0.12272
$ws.Range("J2").Value = -0.1161749333333333
$ws.Range("K2").Value = 1.341
$ws.Range("L2").Value = 0.042912
$ws.Range("M2").Value = 1.42
$ws.Range("N2").Value = 0.02123205741626794
$ws.Range("O2").Value = 1.058911260253542
$ws.Range("P2").Value = 1.42
$ws.Range("Q2").Value = 0.02123205741626794
$ws.Range("R2").Value = 1.058911260253542
$ws.Range("U2").Value = 14.622
$ws.Range("V2").Value = 0.2186303827751196
$ws.Range("W2").Value = 0.002736318407960199
$ws.Range("X2").Value = 0.07531919405177684
$ws.Range("Y2").Value = -0.07258287564381664
$ws.Range("Z2").Value = 2.006807089648086
$ws.Range("AA2").Value = 0.1918965517241379
$ws.Range("AB2").Value = 0.07524554293806116
$ws.Range("AC2").Value = 0.1162014495898479
$ws.Range("AD2").Value = 0.9139999999999999
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 0.9139999999999999
$ws.Range("AG2").Value = -13.708
$ws.Range("AH2").Value = 0.01348201905773372
$ws.Range("AI2").Value = 0.03306323252785415
$ws.Range("AJ2").Value = -0.2578048597005943
$ws.Range("AK2").Value = -1.052680079864844
$ws.Range("AL2").Value = 0.08599999999999999
$ws.Range("AM2").Value = 0.08599999999999999
$ws.Range("AN2").Value = -0.2894236858771374
$ws.Range("AO2").Value = -44.59302325581395
$ws.Range("AP2").Value = 4.340721975934136
$ws.Range("AQ2").Value = -44.59302325581395
$ws.Range("G3").Value = -3.428571428571428
$ws.Range("H3").Value = -3.428571428571428
$ws.Range("I3").Value = -4.32919254658385
$ws.Range("J3").Value = -4.32919254658385
$ws.Range("K3").Value = -2.33
$ws.Range("L3").Value = -1.447204968944099
$ws.Range("M3").Value = -0
$ws.Range("N3").Value = -0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = -0
$ws.Range("Q3").Value = -0
$ws.Range("R3").Value = 0
$ws.Range("U3").Value = 6.57
$ws.Range("V3").Value = 0.4184713375796179
$ws.Range("W3").Value = -0.535632183908046
$ws.Range("X3").Value = 0.07531919405177684
$ws.Range("Y3").Value = -0.6109513779598228
$ws.Range("Z3").Value = -0.8256410256410256
$ws.Range("AA3").Value = 3.574358974358974
$ws.Range("AB3").Value = 0.07524554293806116
$ws.Range("AC3").Value = 3.499113431420913
$ws.Range("AD3").Value = 0.07099999999999999
$ws.Range("AF3").Value = 0.07099999999999999
$ws.Range("AG3").Value = -6.499000000000001
$ws.Range("AH3").Value = 0.004501933929364023
$ws.Range("AI3").Value = 0.05178701677607585
$ws.Range("AJ3").Value = -0.7063362677969788
$ws.Range("AK3").Value = 1.250048086170417
$ws.Range("AN3").Value = -0.01033478893740902
$ws.Range("AP3").Value = 0.9459970887918486
$ws.Range("D4").Value = 0.189
$ws.Range("E4").Value = 0.214
$ws.Range("G4").Value = 0.1733333333333333
$ws.Range("H4").Value = 0.1733333333333333
$ws.Range("I4").Value = 0.1177777777777778
$ws.Range("J4").Value = 0.09893333333333333
$ws.Range("K4").Value = 3.66
$ws.Range("L4").Value = 0.1355555555555556
$ws.Range("M4").Value = 1.42
$ws.Range("N4").Value = 0.03127753303964757
$ws.Range("O4").Value = 0.3879781420765027
$ws.Range("P4").Value = 1.42
$ws.Range("Q4").Value = 0.03127753303964757
$ws.Range("R4").Value = 0.3879781420765027
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 7.77
$ws.Range("V4").Value = 0.1711453744493392
$ws.Range("W4").Value = 0.2140350877192982
$ws.Range("X4").Value = 0.07600253935927474
$ws.Range("Y4").Value = 0.1380325483600235
$ws.Range("Z4").Value = 1.939655172413793
$ws.Range("AA4").Value = 0.1918965517241379
$ws.Range("AB4").Value = 0.07569510213428998
$ws.Range("AC4").Value = 0.1162014495898479
$ws.Range("AD4").Value = 0.834
$ws.Range("AE4").Value = 0
$ws.Range("AF4").Value = 0.834
$ws.Range("AG4").Value = -6.936
$ws.Range("AH4").Value = 0.01803867283817104
$ws.Range("AI4").Value = 0.03668514119820533
$ws.Range("AJ4").Value = -0.180324459234609
$ws.Range("AK4").Value = -0.4635124298315959
$ws.Range("AL4").Value = 0.073
$ws.Range("AM4").Value = 0.073
$ws.Range("AN4").Value = 0.2235924932975871
$ws.Range("AO4").Value = 43.56164383561644
$ws.Range("AP4").Value = -1.859517426273458
$ws.Range("AQ4").Value = 43.56164383561644
$ws.Range("G5").Value = 0.02462121212121212
$ws.Range("H5").Value = 0.02462121212121212
$ws.Range("I5").Value = -0.01704545454545454
$ws.Range("J5").Value = -0.01704545454545454
$ws.Range("K5").Value = 0.011
$ws.Range("L5").Value = 0.004166666666666667
$ws.Range("O5").Value = -0
$ws.Range("R5").Value = -0
$ws.Range("U5").Value = 0.282
$ws.Range("V5").Value = 0.04878892733564013
$ws.Range("W5").Value = 0.002736318407960199
$ws.Range("X5").Value = 0.07517287027675103
$ws.Range("Y5").Value = -0.07243655186879083
$ws.Range("Z5").Value = 0.7329261521377014
$ws.Range("AA5").Value = -0.01249305941143809
$ws.Range("AB5").Value = 0.07514766349911435
$ws.Range("AC5").Value = -0.08764072291055244
$ws.Range("AD5").Value = 0.008999999999999999
$ws.Range("AE5").Value = 0
$ws.Range("AF5").Value = 0.008999999999999999
$ws.Range("AG5").Value = -0.273
$ws.Range("AH5").Value = 0.001554672655035412
$ws.Range("AI5").Value = 0.002543091268719978
$ws.Range("AJ5").Value = -0.04957327038314871
$ws.Range("AK5").Value = -0.08381946576604236
$ws.Range("AL5").Value = 0.013
$ws.Range("AM5").Value = 0.013
$ws.Range("AN5").Value = -0.5
$ws.Range("AO5").Value = -3.461538461538462
$ws.Range("AP5").Value = 15.16666666666667
$ws.Range("AQ5").Value = -3.461538461538462
